# edit.ps1 - applies the "Cyber Security" -> "Chemistry" rewrite described by the diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Builds a single <w:r> run.
#   $text   - text content (may be empty when $breakOnly is used)
#   $sz     - font size (half-points) as string, or $null to omit <w:sz>
#   $breakBefore - if $true, emit a <w:br/> before the text
#   $breakOnly   - if $true, the run is just a manual line break (<w:br/>), $text is ignored
#   $pageBreakBefore - if $true, emit a <w:lastRenderedPageBreak/> before the text
function New-RunXml {
    param(
        [string]$text,
        $sz = $null,
        [bool]$breakBefore = $false,
        [bool]$breakOnly = $false,
        [bool]$pageBreakBefore = $false
    )

    $rpr = "<w:rPr><w:rFonts w:ascii=`"Times New Roman`" w:hAnsi=`"Times New Roman`"/><w:color w:val=`"000000`"/>"
    if ($sz) {
        $rpr += "<w:sz w:val=`"$sz`"/>"
    }
    $rpr += "</w:rPr>"

    $body = ""
    if ($breakOnly) {
        $body = "<w:br/>"
    } else {
        if ($pageBreakBefore) {
            $body += "<w:lastRenderedPageBreak/>"
        }
        if ($breakBefore) {
            $body += "<w:br/>"
        }
        $needsPreserve = ($text.Length -gt 0) -and (($text.StartsWith(" ")) -or ($text.EndsWith(" ")))
        if ($needsPreserve) {
            $body += "<w:t xml:space=`"preserve`">$text</w:t>"
        } else {
            $body += "<w:t>$text</w:t>"
        }
    }
    return "<w:r>$rpr$body</w:r>"
}

# Wraps a sequence of run-xml fragments into a full WordOpenXML package fragment
# suitable for Range.InsertXML, representing a single paragraph's contents.
function New-ParagraphPackageXml {
    param([string]$runsXml)

    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Replaces the run-content of paragraph number $paraIndex (1-based) with the
# supplied run-xml, leaving the paragraph mark (and therefore its w:pPr) intact.
function Set-ParagraphRuns {
    param([int]$paraIndex, [string]$runsXml)

    $para = $d.Paragraphs($paraIndex)
    $full = $para.Range
    $r = $d.Range($full.Start, $full.End - 1)
    $pkg = New-ParagraphPackageXml $runsXml
    $r.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1) Title
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Cyber Security: The Digital Shield", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Unveiling the Marvels of Chemistry: A Journey of Transformation and Discovery", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Author name paragraph: "Andrew Davenport" -> "Dr" + "." + " Erika Martinez"
# ---------------------------------------------------------------------------
$authorRuns = (New-RunXml "Dr" "36") + (New-RunXml "." "36") + (New-RunXml " Erika Martinez" "36")
Set-ParagraphRuns 2 $authorRuns

# ---------------------------------------------------------------------------
# 3) Email paragraph: "andrew" "." "davenport@blackrock" "." "net"
#    -> "emartinez@highschooledu" "." "org"
# ---------------------------------------------------------------------------
$emailRuns = (New-RunXml "emartinez@highschooledu" "32") + (New-RunXml "." "32") + (New-RunXml "org" "32")
Set-ParagraphRuns 3 $emailRuns

# ---------------------------------------------------------------------------
# 5) Main body paragraph - full rewrite (chemistry essay)
# ---------------------------------------------------------------------------
$bodyRuns = ""
$bodyRuns += New-RunXml "In the vast tapestry of human knowledge, chemistry stands as a vibrant and transformative thread, weaving together the elements of the universe to reveal the wonders of matter and its intricate interactions" "24"
$bodyRuns += New-RunXml "." "24"
$bodyRuns += New-RunXml " It is the language of nature's symphony, guiding us through the profound mysteries of chemical reactions, the compositions of substances, and the fundamental principles that govern the behaviour of matter" "24"
$bodyRuns += New-RunXml "." "24"
$bodyRuns += New-RunXml " As we embark on this enthralling adventure into the world of chemistry, let us unravel the enigmatic secrets of this fascinating science, unveiling its profound impact on our lives" "24"
$bodyRuns += New-RunXml "." "24"
$bodyRuns += New-RunXml "" "24" $false $true
$bodyRuns += New-RunXml "In the realm of chemistry, we witness the ceaseless dance of atoms and molecules, orchestrated by the enigmatic forces of attraction and repulsion" "24" $true
$bodyRuns += New-RunXml "." "24"
$bodyRuns += New-RunXml " These fundamental particles engage in a ceaseless choreography, forming compounds and rearranging their structures, giving rise to the mesmerizing diversity of substances that make up our world" "24"
$bodyRuns += New-RunXml "." "24"
$bodyRuns += New-RunXml " From the air we breathe and the water we drink, to the food we consume and the medicines that heal us, chemistry lies at the heart of all life" "24"
$bodyRuns += New-RunXml "." "24"
$bodyRuns += New-RunXml " It is the driving force behind the symphony of life, governing the intricate interactions between organisms in the intricate web of ecosystems" "24"
$bodyRuns += New-RunXml "." "24"
$bodyRuns += New-RunXml "" "24" $false $true
$bodyRuns += New-RunXml "Through the lens of chemistry, we can unravel the complexities of chemical reactions, understanding how substances transform from one state to another, releasing energy or undergoing profound changes in their properties" "24" $true
$bodyRuns += New-RunXml "." "24"
$bodyRuns += New-RunXml " We uncover the secrets of catalysts, molecules that accelerate these transformations, allowing us to harness nature's power to create new substances and materials" "24"
$bodyRuns += New-RunXml "." "24"
$bodyRuns += New-RunXml " Chemistry empowers us with the ability to synthesize drugs that combat diseases, develop materials with extraordinary properties, and create sustainable energy sources, all of which have the potential to shape a better future for humankind" "24"
$bodyRuns += New-RunXml "." "24"
Set-ParagraphRuns 5 $bodyRuns

# ---------------------------------------------------------------------------
# 6) "Summary" heading - text unchanged, just ensure font later via global pass
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 7) Final summary paragraph - full rewrite
# ---------------------------------------------------------------------------
$finalRuns = ""
$finalRuns += New-RunXml "Chemistry is the vibrant language of nature, revealing the secrets of matter and its intricate interactions" $null
$finalRuns += New-RunXml "." $null
$finalRuns += New-RunXml " From the symphony of atomic dances to the transformative power of chemical reactions, chemistry plays a pivotal role in shaping our world" $null
$finalRuns += New-RunXml "." $null
$finalRuns += New-RunXml " It holds the key to understanding the composition and behaviour of substances, leading to advancements in medicine, technology, " $null
$finalRuns += New-RunXml "and sustainable energy" $null $false $false $true
$finalRuns += New-RunXml "." $null
$finalRuns += New-RunXml " Chemistry empowers us to comprehend and manipulate the world around us, fostering progress and innovation while deepening our appreciation for the marvels of nature" $null
$finalRuns += New-RunXml "." $null
Set-ParagraphRuns 7 $finalRuns

# ---------------------------------------------------------------------------
# 8) Add a new empty paragraph at the very end of the document (before sectPr)
# ---------------------------------------------------------------------------
$endRange = $d.Paragraphs(7).Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 9) Global font fix: "TimesNewToman" -> "Times New Roman" for the runs that
#    were not rebuilt via InsertXML above (the title and the "Summary"
#    heading). Font names live in run properties, not document text, so a
#    text Find/Replace cannot reach them - set Font.Name on the run ranges
#    directly instead (excluding the trailing paragraph mark so the
#    paragraph-mark's own rPr is left untouched).
# ---------------------------------------------------------------------------
function Set-ParagraphFont {
    param([int]$paraIndex, [string]$fontName)

    $full = $d.Paragraphs($paraIndex).Range
    $r = $d.Range($full.Start, $full.End - 1)
    $r.Font.Name = $fontName
}

Set-ParagraphFont 1 "Times New Roman"
Set-ParagraphFont 6 "Times New Roman"
